$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 4 (C4:G4) down to row 5 (C5:G5) so the new
# cells inherit the same number formats / styles as the row above.
$ws.Range("C4:G4").Copy()
$ws.Range("C5:G5").PasteSpecial(-4122)

# Populate the new values for row 5
$ws.Range("C5").Value = 50.445
$ws.Range("D5").Value = 75.66
$ws.Range("E5").Value = "E522 BP1 9V ALCALINA ENERGIZER MAX"
$ws.Range("F5").Value = 0.21
$ws.Range("G5").Value = 0.21

# Update the active selection to match the saved workbook state
$ws.Activate()
$ws.Range("E12").Select()
